# ms revision on 20220828-2
#
# Applies:
#   - TextBox 76: "AT3G13920.1" -> "AT3G13910.2"
#   - TextBox 14 ("(a)"): reposition/resize + text "(a)" -> "A"
#   - TextBox 19 ("(b)"): resize + text "(b)" -> "B"
#   - TextBox 24 ("(c)"): resize + text "(c)" -> "C"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> points helper (914400 EMU per inch, 72 points per inch).
# The COM layer stores Left/Top/Width/Height as single-precision point
# values and truncates when converting back to EMU on save, so nudge by
# half an EMU to land exactly on the target EMU after that round trip.
function EmuToPt($emu) {
    return ($emu + 0.5) / 914400 * 72
}

# --- Gene id label correction -------------------------------------------
$geneShape = $s.Shapes.Item("TextBox 76")
$geneShape.TextFrame.TextRange.Text = "AT3G13910.2"

# --- Panel label (a) -> A, moved down and narrowed -----------------------
$panelA = $s.Shapes.Item("TextBox 14")
$panelA.Left = EmuToPt 1575331
$panelA.Top = EmuToPt 1344157
$panelA.Width = EmuToPt 338554
$panelA.Height = EmuToPt 369332
$panelA.TextFrame.TextRange.Text = "A"

# --- Panel label (b) -> B, narrowed ---------------------------------------
$panelB = $s.Shapes.Item("TextBox 19")
$panelB.Left = EmuToPt 1583139
$panelB.Top = EmuToPt 3616090
$panelB.Width = EmuToPt 338554
$panelB.Height = EmuToPt 369332
$panelB.TextFrame.TextRange.Text = "B"

# --- Panel label (c) -> C, narrowed ---------------------------------------
$panelC = $s.Shapes.Item("TextBox 24")
$panelC.Left = EmuToPt 8485078
$panelC.Top = EmuToPt 3645415
$panelC.Width = EmuToPt 351378
$panelC.Height = EmuToPt 369332
$panelC.TextFrame.TextRange.Text = "C"
